$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 3520.6
$ws.Range("I31").Value = 201
$ws.Range("J31").Value = 8500
$ws.Range("K31").Value = 603
$ws.Range("L31").Value = 25500
$ws.Range("M31").Value = -373

$ws.Range("H39").Value = 868.5714
$ws.Range("I39").Value = 136.16667
$ws.Range("J39").Value = 5263
$ws.Range("K39").Value = 408.50001
$ws.Range("L39").Value = 15789
$ws.Range("M39").Value = -112.50001
$ws.Range("N39").Value = -16381

$ws.Range("H41").Value = 2009.5
$ws.Range("I41").Value = 1246.6666
$ws.Range("J41").Value = 4298
$ws.Range("K41").Value = 1246.6666
$ws.Range("L41").Value = 4298
$ws.Range("M41").Value = -806.6666
$ws.Range("N41").Value = -5178

$ws.Range("H80").Value = 828.2857
$ws.Range("I80").Value = 469
$ws.Range("J80").Value = 940.5625
$ws.Range("K80").Value = 1407
$ws.Range("L80").Value = 2821.6875
$ws.Range("M80").Value = -409
$ws.Range("N80").Value = -4817.6875

$ws.Range("H83").Value = 828.2857
$ws.Range("I83").Value = 469
$ws.Range("J83").Value = 940.5625
$ws.Range("K83").Value = 4221
$ws.Range("L83").Value = 8465.0625
$ws.Range("M83").Value = 771
$ws.Range("N83").Value = -18449.0625

$ws.Range("H99").Value = 1872.2
$ws.Range("I99").Value = 565.3333
$ws.Range("J99").Value = 3832.5
$ws.Range("K99").Value = 1695.9999
$ws.Range("L99").Value = 11497.5
$ws.Range("M99").Value = -197.9999

$ws.Range("H105").Value = 35268.4
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 35268.4
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 35268.4
$ws.Range("N105").Value = -42256.4

$ws.Range("H106").Value = 7500
$ws.Range("I106").Value = 7500
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 7500
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -6869
$ws.Range("N106").ClearContents()

$ws.Range("H117").Value = 70433
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 70433
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 70433
$ws.Range("N117").Value = -79611

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2906.6365
$ws.Range("I45").Value = 2121.75
$ws.Range("J45").Value = 4999.6665
$ws.Range("K45").Value = 2121.75
$ws.Range("L45").Value = 4999.6665
$ws.Range("M45").Value = -1744.75

$ws.Range("H61").Value = 4079.6155
$ws.Range("I61").Value = 3639.7273
$ws.Range("J61").Value = 6499
$ws.Range("K61").Value = 3639.7273
$ws.Range("L61").Value = 6499
$ws.Range("M61").Value = -3427.7273

$ws.Range("H122").Value = 2298.8
$ws.Range("I122").Value = 1998.6666
$ws.Range("J122").Value = 2749
$ws.Range("K122").Value = 5995.9998
$ws.Range("L122").Value = 8247
$ws.Range("M122").Value = -3545.9998
$ws.Range("N122").Value = -13147

$ws.Range("H132").Value = 3715.8462
$ws.Range("I132").Value = 1811.8889
$ws.Range("J132").Value = 7999.75
$ws.Range("K132").Value = 5435.6667
$ws.Range("L132").Value = 23999.25
$ws.Range("M132").Value = -2905.6667

$ws.Range("H136").Value = 4079.6155
$ws.Range("I136").Value = 3639.7273
$ws.Range("J136").Value = 6499
$ws.Range("K136").Value = 10919.1819
$ws.Range("L136").Value = 19497
$ws.Range("M136").Value = -8369.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1562.4736
$ws.Range("I105").Value = 1246.7693
$ws.Range("J105").Value = 2246.5
$ws.Range("K105").Value = 1246.7693
$ws.Range("L105").Value = 2246.5
$ws.Range("M105").Value = 500.2307000000001

$ws.Range("H107").Value = 5025.304
$ws.Range("I107").Value = 1144.091
$ws.Range("J107").Value = 8583.083000000001
$ws.Range("K107").Value = 1144.091
$ws.Range("L107").Value = 8583.083000000001
$ws.Range("M107").Value = 775.9090000000001

$ws.Range("H132").Value = 90000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 90000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 90000
$ws.Range("N132").Value = -100120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H105").Value = 1320.5555
$ws.Range("I105").Value = 1320.5555
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1320.5555
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 426.4445000000001

$ws.Range("H106").Value = 35780
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 35780
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 35780
$ws.Range("N106").Value = -38304

$ws.Range("H107").Value = 341.83334
$ws.Range("I107").Value = 652.5
$ws.Range("J107").Value = 186.5
$ws.Range("K107").Value = 652.5
$ws.Range("L107").Value = 186.5
$ws.Range("M107").Value = 1267.5
$ws.Range("N107").Value = -4026.5

$ws.Range("H132").Value = 5203.154
$ws.Range("I132").Value = 4455.875
$ws.Range("J132").Value = 6398.8
$ws.Range("K132").Value = 13367.625
$ws.Range("L132").Value = 19196.4
$ws.Range("M132").Value = -10837.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1474.3334
$ws.Range("I58").Value = 1036
$ws.Range("J58").Value = 2022.25
$ws.Range("K58").Value = 3108
$ws.Range("L58").Value = 6066.75
$ws.Range("M58").Value = -2980
$ws.Range("N58").Value = -6322.75

$ws.Range("H115").Value = 1467.5
$ws.Range("I115").Value = 1290
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 3870
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -2695
$ws.Range("N115").Value = -8350

$ws.Range("H137").Value = 2932
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 3898.5
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 11695.5
$ws.Range("M137").Value = 2103
$ws.Range("N137").Value = -21895.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 27443
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 27443
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 27443
$ws.Range("N47").Value = -28579
$ws.Range("M47").ClearContents()

$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3500
$ws.Range("N80").Value = -5496

$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17500
$ws.Range("N83").Value = -27484

$ws.Range("H97").Value = 775.5
$ws.Range("I97").Value = 811.2308
$ws.Range("J97").Value = 311
$ws.Range("K97").Value = 811.2308
$ws.Range("L97").Value = 311
$ws.Range("M97").Value = -315.2308
$ws.Range("N97").Value = -1303

$ws.Range("H108").Value = 74499.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 74499.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 74499.5
$ws.Range("N108").Value = -82179.5

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 54629.45
$ws.Range("I132").Value = 59588.332
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 178764.996
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -176234.996
$ws.Range("N132").Value = -35058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2595.1428
$ws.Range("I55").Value = 1833.4
$ws.Range("J55").Value = 4499.5
$ws.Range("K55").Value = 1833.4
$ws.Range("L55").Value = 4499.5
$ws.Range("M55").Value = -1660.4

$ws.Range("H61").Value = 6680.5454
$ws.Range("I61").Value = 5581
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 5581
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -5379

$ws.Range("H100").Value = 7264.5186
$ws.Range("I100").Value = 5766.5713
$ws.Range("J100").Value = 7788.8
$ws.Range("K100").Value = 5766.5713
$ws.Range("L100").Value = 7788.8
$ws.Range("M100").Value = -5225.5713

$ws.Range("H113").Value = 6680.5454
$ws.Range("I113").Value = 5581
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 5581
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -3411

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3501
$ws.Range("I96").Value = 2751.5
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 2751.5
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -1378.5

$ws.Range("H113").Value = 837.38464
$ws.Range("I113").Value = 1193.1428
$ws.Range("J113").Value = 422.33334
$ws.Range("K113").Value = 3579.4284
$ws.Range("L113").Value = 1267.00002
$ws.Range("M113").Value = -1409.4284
$ws.Range("N113").Value = -5607.000019999999

$ws.Range("H132").Value = 1944.6
$ws.Range("I132").Value = 1944.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5833.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3303.799999999999
